$wb = $excel.ActiveWorkbook

# Overview sheet: bump pending "Ready for handoff" / "Handback transform failed"
# handoff-date column (D) to the refreshed timestamp.
$ws1 = $wb.Worksheets.Item("Overview")
foreach ($r in 7,10,11,12,13,14,15,16) {
    $ws1.Cells.Item($r, 4).Value = "2016-22-14 03:22:49"
}

# zh-cn sheet: bump "Latest Handoff Datetime" column (E) for the same rows.
$ws2 = $wb.Worksheets.Item("zh-cn")
foreach ($r in 7,10,11,12,13,14,15,16) {
    $ws2.Cells.Item($r, 5).Value = "2016-03-14 03:22:45"
}

# de-de sheet: bump "Latest Handoff Datetime" column (E) for the same rows.
$ws3 = $wb.Worksheets.Item("de-de")
foreach ($r in 7,10,11,12,13,14,15,16) {
    $ws3.Cells.Item($r, 5).Value = "2016-03-14 03:22:49"
}
